# Remove the duplicate "Resumo Executivo / Diagnostico Estrutural / Lacunas /
# Classificacao / Observacoes XALQ" block that was left over after the page
# break near the end of the document, and drop the now-orphaned page break
# itself (the paragraph that held it stays, but becomes empty).

$d = $word.ActiveDocument

# Locate the paragraph that holds the manual page break ("\f" i.e. chr(12))
# near the tail of the document, then find where the duplicated block (the
# one starting with "Resumo Executivo") starts right after it.
$breakParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.Length -ge 1 -and [int][char]$t[0] -eq 12) {
        $breakParaIndex = $i
    }
}

if ($breakParaIndex -gt 0 -and $breakParaIndex -lt $d.Paragraphs.Count) {
    # Delete everything from the paragraph right after the page-break
    # paragraph through the end of the document body content (the stray
    # duplicated "Resumo Executivo ... Observacoes XALQ" section).
    $dupStart = $d.Paragraphs($breakParaIndex + 1).Range.Start
    $docEnd = $d.Content.End
    $dupRange = $d.Range($dupStart, $docEnd)
    $dupRange.Delete()

    # Now strip just the page-break character itself, leaving the paragraph
    # mark (and its pPr/formatting) intact as an empty paragraph.
    $breakPara = $d.Paragraphs($breakParaIndex)
    $breakCharStart = $breakPara.Range.Start
    $breakCharRange = $d.Range($breakCharStart, $breakCharStart + 1)
    $breakCharRange.Delete()
}
